# Insert two new weekly price rows (Poroto granado, Femacal de La Calera)
# right after the existing row 145, pushing all subsequent rows down by two
# positions (old row 146 -> new row 148, ..., old row 242 -> new row 244).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 146 downward by inserting two blank rows at 146:147.
$ws.Rows("146:147").Insert()

# Fill in the first new row (146).
$ws.Range("A146").Value = 3
$ws.Range("B146").Value = "Femacal de La Calera"
$ws.Range("C146").Value = "Coquimbo"
$ws.Range("D146").Value = 44960
$ws.Range("E146").Value = 5
$ws.Range("F146").Value = 100112030
$ws.Range("G146").Value = "Poroto granado"
$ws.Range("H146").Value = "Sin especificar"
$ws.Range("I146").Value = "Primera"
$ws.Range("J146").Value = 60
$ws.Range("K146").Value = 43000
$ws.Range("L146").Value = 44000
$ws.Range("M146").Value = 43500
$ws.Range("N146").Value = "$/saco 25 kilos"
$ws.Range("O146").Value = "Provincia de Quillota"
$ws.Range("P146").Value = 1740
$ws.Range("Q146").Value = 25
$ws.Range("R146").Value = "Hortaliza"

# Fill in the second new row (147).
$ws.Range("A147").Value = 3
$ws.Range("B147").Value = "Femacal de La Calera"
$ws.Range("C147").Value = "Coquimbo"
$ws.Range("D147").Value = 44960
$ws.Range("E147").Value = 5
$ws.Range("F147").Value = 100112030
$ws.Range("G147").Value = "Poroto granado"
$ws.Range("H147").Value = "Sin especificar"
$ws.Range("I147").Value = "Segunda"
$ws.Range("J147").Value = 30
$ws.Range("K147").Value = 39000
$ws.Range("L147").Value = 39000
$ws.Range("M147").Value = 39000
$ws.Range("N147").Value = "$/saco 25 kilos"
$ws.Range("O147").Value = "Provincia de Quillota"
$ws.Range("P147").Value = 1560
$ws.Range("Q147").Value = 25
$ws.Range("R147").Value = "Hortaliza"
